$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = '@'
$cell.Value = '68.133.49'
$cell.Style = 'Normal'
$cell = $ws.Range("E2")
$cell.NumberFormat = '@'
$cell.Value = '  +2.74%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D3")
$cell.NumberFormat = '@'
$cell.Value = '3.809.72'
$cell.Style = 'Normal'
$cell = $ws.Range("E3")
$cell.NumberFormat = '@'
$cell.Value = '  +7.04%  '
$cell.Style = 'Normal'
$cell = $ws.Range("E4")
$cell.NumberFormat = '@'
$cell.Value = '  +0.28%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D5")
$cell.NumberFormat = '@'
$cell.Value = '418.25'
$cell.Style = 'Normal'
$cell = $ws.Range("E5")
$cell.NumberFormat = '@'
$cell.Value = '  -0.12%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D6")
$cell.NumberFormat = '@'
$cell.Value = '138.73'
$cell.Style = 'Normal'
$cell = $ws.Range("E6")
$cell.NumberFormat = '@'
$cell.Value = '  +5.01%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D7")
$cell.NumberFormat = '@'
$cell.Value = '3.795.64'
$cell.Style = 'Normal'
$cell = $ws.Range("E7")
$cell.NumberFormat = '@'
$cell.Value = '  +6.95%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D8")
$cell.NumberFormat = '@'
$cell.Value = '0.650'
$cell.Style = 'Normal'
$cell = $ws.Range("E8")
$cell.NumberFormat = '@'
$cell.Value = '  -1.56%  '
$cell.Style = 'Normal'
$cell = $ws.Range("E9")
$cell.NumberFormat = '@'
$cell.Value = '  -0.02%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D10")
$cell.NumberFormat = '@'
$cell.Value = '0.772'
$cell.Style = 'Normal'
$cell = $ws.Range("E10")
$cell.NumberFormat = '@'
$cell.Value = '  -1.66%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D11")
$cell.NumberFormat = '@'
$cell.Value = '0.183'
$cell.Style = 'Normal'
$cell = $ws.Range("E11")
$cell.NumberFormat = '@'
$cell.Value = '  +8.10%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D12")
$cell.NumberFormat = '@'
$cell.Value = '0.0000398'
$cell.Style = 'Normal'
$cell = $ws.Range("E12")
$cell.NumberFormat = '@'
$cell.Value = '  +38.43%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D13")
$cell.NumberFormat = '@'
$cell.Value = '43.16'
$cell.Style = 'Normal'
$cell = $ws.Range("E13")
$cell.NumberFormat = '@'
$cell.Value = '  -0.55%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D14")
$cell.NumberFormat = '@'
$cell.Value = '10.40'
$cell.Style = 'Normal'
$cell = $ws.Range("D15")
$cell.NumberFormat = '@'
$cell.Value = '4.396.83'
$cell.Style = 'Normal'
$cell = $ws.Range("E15")
$cell.NumberFormat = '@'
$cell.Value = '  +6.71%  '
$cell.Style = 'Normal'
$cell = $ws.Range("E16")
$cell.NumberFormat = '@'
$cell.Value = '  -0.54%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D17")
$cell.NumberFormat = '@'
$cell.Value = '3.816.51'
$cell.Style = 'Normal'
$cell = $ws.Range("E17")
$cell.NumberFormat = '@'
$cell.Value = '  +6.55%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D18")
$cell.NumberFormat = '@'
$cell.Value = '20.64'
$cell.Style = 'Normal'
$cell = $ws.Range("E18")
$cell.NumberFormat = '@'
$cell.Value = '  +0.09%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D19")
$cell.NumberFormat = '@'
$cell.Value = '13.21'
$cell.Style = 'Normal'
$cell = $ws.Range("E19")
$cell.NumberFormat = '@'
$cell.Value = '  +3.62%  '
$cell.Style = 'Normal'
$cell = $ws.Range("E20")
$cell.NumberFormat = '@'
$cell.Value = '  +1.76%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D21")
$cell.NumberFormat = '@'
$cell.Value = '68.263.58'
$cell.Style = 'Normal'
$cell = $ws.Range("E21")
$cell.NumberFormat = '@'
$cell.Value = '  +3.01%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D22")
$cell.NumberFormat = '@'
$cell.Value = '443.22'
$cell.Style = 'Normal'
$cell = $ws.Range("E22")
$cell.NumberFormat = '@'
$cell.Value = '  -2.02%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D23")
$cell.NumberFormat = '@'
$cell.Value = '15.37'
$cell.Style = 'Normal'
$cell = $ws.Range("E23")
$cell.NumberFormat = '@'
$cell.Value = '  +16.41%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D24")
$cell.NumberFormat = '@'
$cell.Value = '89.80'
$cell.Style = 'Normal'
$cell = $ws.Range("E24")
$cell.NumberFormat = '@'
$cell.Value = '  -0.70%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D25")
$cell.NumberFormat = '@'
$cell.Value = '3.11'
$cell.Style = 'Normal'
$cell = $ws.Range("E25")
$cell.NumberFormat = '@'
$cell.Value = '  -4.56%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D26")
$cell.NumberFormat = '@'
$cell.Value = '38.12'
$cell.Style = 'Normal'
$cell = $ws.Range("E26")
$cell.NumberFormat = '@'
$cell.Value = '  +10.71%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D27")
$cell.NumberFormat = '@'
$cell.Value = '3.31'
$cell.Style = 'Normal'
$cell = $ws.Range("E27")
$cell.NumberFormat = '@'
$cell.Value = '  -2.37%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D28")
$cell.NumberFormat = '@'
$cell.Value = '9.91'
$cell.Style = 'Normal'
$cell = $ws.Range("E28")
$cell.NumberFormat = '@'
$cell.Value = '  -1.41%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D29")
$cell.NumberFormat = '@'
$cell.Value = '5.16'
$cell.Style = 'Normal'
$cell = $ws.Range("E29")
$cell.NumberFormat = '@'
$cell.Value = '  +6.62%  '
$cell.Style = 'Normal'
$cell = $ws.Range("B30")
$cell.NumberFormat = '@'
$cell.Value = 'Cosmos'
$cell.Style = 'Normal'
$cell = $ws.Range("C30")
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell.Style = 'Normal'
$cell = $ws.Range("D30")
$cell.NumberFormat = '@'
$cell.Value = '12.65'
$cell.Style = 'Normal'
$cell = $ws.Range("E30")
$cell.NumberFormat = '@'
$cell.Value = '  +1.47%  '
$cell.Style = 'Normal'
$cell = $ws.Range("B31")
$cell.NumberFormat = '@'
$cell.Value = 'Hedera'
$cell.Style = 'Normal'
$cell = $ws.Range("C31")
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$cell.Style = 'Normal'
$cell = $ws.Range("D31")
$cell.NumberFormat = '@'
$cell.Value = '0.123'
$cell.Style = 'Normal'
$cell = $ws.Range("E31")
$cell.NumberFormat = '@'
$cell.Value = '  +5.00%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D32")
$cell.NumberFormat = '@'
$cell.Value = '2.76'
$cell.Style = 'Normal'
$cell = $ws.Range("E32")
$cell.NumberFormat = '@'
$cell.Value = '  -0.77%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D33")
$cell.NumberFormat = '@'
$cell.Value = '7.15'
$cell.Style = 'Normal'
$cell = $ws.Range("E33")
$cell.NumberFormat = '@'
$cell.Value = '  -2.35%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D34")
$cell.NumberFormat = '@'
$cell.Value = '41.64'
$cell.Style = 'Normal'
$cell = $ws.Range("E34")
$cell.NumberFormat = '@'
$cell.Value = '  +6.07%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D35")
$cell.NumberFormat = '@'
$cell.Value = '0.162'
$cell.Style = 'Normal'
$cell = $ws.Range("E35")
$cell.NumberFormat = '@'
$cell.Value = '  +0.14%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D36")
$cell.NumberFormat = '@'
$cell.Value = '57.84'
$cell.Style = 'Normal'
$cell = $ws.Range("E36")
$cell.NumberFormat = '@'
$cell.Value = '  -0.03%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D37")
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$cell = $ws.Range("E37")
$cell.NumberFormat = '@'
$cell.Value = '  +0.09%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D38")
$cell.NumberFormat = '@'
$cell.Value = '0.0490'
$cell.Style = 'Normal'
$cell = $ws.Range("E38")
$cell.NumberFormat = '@'
$cell.Value = '  -3.38%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D39")
$cell.NumberFormat = '@'
$cell.Value = '3.04'
$cell.Style = 'Normal'
$cell = $ws.Range("E39")
$cell.NumberFormat = '@'
$cell.Value = '  +30.15%  '
$cell.Style = 'Normal'
$cell = $ws.Range("B40")
$cell.NumberFormat = '@'
$cell.Value = 'Stellar'
$cell.Style = 'Normal'
$cell = $ws.Range("C40")
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$cell.Style = 'Normal'
$cell = $ws.Range("D40")
$cell.NumberFormat = '@'
$cell.Value = '0.148'
$cell.Style = 'Normal'
$cell = $ws.Range("E40")
$cell.NumberFormat = '@'
$cell.Value = '  -0.36%  '
$cell.Style = 'Normal'
$cell = $ws.Range("B41")
$cell.NumberFormat = '@'
$cell.Value = 'PEPE'
$cell.Style = 'Normal'
$cell = $ws.Range("C41")
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$cell.Style = 'Normal'
$cell = $ws.Range("D41")
$cell.NumberFormat = '@'
$cell.Value = '0.0₃0698'
$cell.Style = 'Normal'
$cell = $ws.Range("E41")
$cell.NumberFormat = '@'
$cell.Value = '  -5.69%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D42")
$cell.NumberFormat = '@'
$cell.Value = '0.998'
$cell.Style = 'Normal'
$cell = $ws.Range("E42")
$cell.NumberFormat = '@'
$cell.Value = '  +0.02%  '
$cell.Style = 'Normal'
$cell = $ws.Range("B43")
$cell.NumberFormat = '@'
$cell.Value = 'EnergySwap'
$cell.Style = 'Normal'
$cell = $ws.Range("C43")
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell.Style = 'Normal'
$cell = $ws.Range("D43")
$cell.NumberFormat = '@'
$cell.Value = '27.77'
$cell.Style = 'Normal'
$cell = $ws.Range("E43")
$cell.NumberFormat = '@'
$cell.Value = '  +28.74%  '
$cell.Style = 'Normal'
$cell = $ws.Range("B44")
$cell.NumberFormat = '@'
$cell.Value = 'ApeXProtocol'
$cell.Style = 'Normal'
$cell = $ws.Range("C44")
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$cell.Style = 'Normal'
$cell = $ws.Range("D44")
$cell.NumberFormat = '@'
$cell.Value = '3.24'
$cell.Style = 'Normal'
$cell = $ws.Range("E44")
$cell.NumberFormat = '@'
$cell.Value = '  +25.18%  '
$cell.Style = 'Normal'
$cell = $ws.Range("B45")
$cell.NumberFormat = '@'
$cell.Value = 'LidoDAOToken'
$cell.Style = 'Normal'
$cell = $ws.Range("C45")
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$cell.Style = 'Normal'
$cell = $ws.Range("D45")
$cell.NumberFormat = '@'
$cell.Value = '3.40'
$cell.Style = 'Normal'
$cell = $ws.Range("E45")
$cell.NumberFormat = '@'
$cell.Value = '  +3.66%  '
$cell.Style = 'Normal'
$cell = $ws.Range("B46")
$cell.NumberFormat = '@'
$cell.Value = 'Monero'
$cell.Style = 'Normal'
$cell = $ws.Range("C46")
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell.Style = 'Normal'
$cell = $ws.Range("D46")
$cell.NumberFormat = '@'
$cell.Value = '148.32'
$cell.Style = 'Normal'
$cell = $ws.Range("E46")
$cell.NumberFormat = '@'
$cell.Value = '  +0.15%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D47")
$cell.NumberFormat = '@'
$cell.Value = '2.10'
$cell.Style = 'Normal'
$cell = $ws.Range("E47")
$cell.NumberFormat = '@'
$cell.Value = '  +4.49%  '
$cell.Style = 'Normal'
$cell = $ws.Range("E48")
$cell.NumberFormat = '@'
$cell.Value = '  -6.15%  '
$cell.Style = 'Normal'
$cell = $ws.Range("B49")
$cell.NumberFormat = '@'
$cell.Value = 'NEARProtocol'
$cell.Style = 'Normal'
$cell = $ws.Range("C49")
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$cell.Style = 'Normal'
$cell = $ws.Range("D49")
$cell.NumberFormat = '@'
$cell.Value = '4.30'
$cell.Style = 'Normal'
$cell = $ws.Range("E49")
$cell.NumberFormat = '@'
$cell.Value = '  -3.75%  '
$cell.Style = 'Normal'
$cell = $ws.Range("B50")
$cell.NumberFormat = '@'
$cell.Value = 'WEMIXToken'
$cell.Style = 'Normal'
$cell = $ws.Range("C50")
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$cell.Style = 'Normal'
$cell = $ws.Range("D50")
$cell.NumberFormat = '@'
$cell.Value = '2.61'
$cell.Style = 'Normal'
$cell = $ws.Range("E50")
$cell.NumberFormat = '@'
$cell.Value = '  -5.99%  '
$cell.Style = 'Normal'
$cell = $ws.Range("D51")
$cell.NumberFormat = '@'
$cell.Value = '0.304'
$cell.Style = 'Normal'
$cell = $ws.Range("E51")
$cell.NumberFormat = '@'
$cell.Value = '  -2.62%  '
$cell.Style = 'Normal'
